$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Change the repayment strategy value scenario from "RBI (India)" to
# "Overdue/Due Fee/Int,Principal" (row 17 -> A17 "repaymentstrategy")
$wsInput.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Update the active selection to reflect the edited cell
$wsInput.Range("B17").Select()
